# Appends a "Things to add" TODO list to the end of the document.
#
# We build the new content as a WordprocessingML fragment and insert it in
# one shot via Range.InsertXML (wrapped in the standard Flat-OPC
# <pkg:package> envelope that Word's InsertXML expects). Doing it this way
# -- rather than via repeated InsertParagraphAfter/TypeText calls -- lets us
# produce exactly the paragraph/run/proofErr structure we want (including
# genuinely empty paragraphs with no runs, and paragraphs that don't inherit
# the preceding list-paragraph numbering/style).

$d = $word.ActiveDocument

# Collapse a range to the very end of the document body (after the last
# existing paragraph, before the end-of-story mark).
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)

$newBodyFragment = (
  '<w:p/>' +
  '<w:p/>' +
  '<w:p/>' +
  '<w:p>' +
    '<w:r><w:t xml:space="preserve">Things to </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>add</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>' +
  '<w:p>' +
    '<w:r><w:t>Cache</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> centroid calculations</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:r><w:t xml:space="preserve">Add proper names for folders via open ai </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>database</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>' +
  '<w:p><w:r><w:t>Progress bar</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t>Error handling</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t>Requirements.txt</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t>Readme</w:t></w:r></w:p>' +
  '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Github</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>' +
  '<w:p/>'
)

$flatOpcXml = (
  '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
      'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $newBodyFragment + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
)

$null = $insertionPoint.InsertXML($flatOpcXml)

Write-Host ("Paragraphs.Count=" + $d.Paragraphs.Count)
